# Insert a new "Optimal" column before the existing "Description" column (E),
# pushing the old column E ("Description") to F, then populate the new
# column with More/Less guidance for each metric row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new column at E; existing E (Description) and its data shift to F.
$ws.Range("E1").EntireColumn.Insert()

# Header for the new column.
$ws.Range("E1").Value = "Optimal"

# Per-row guidance values for the new "Optimal" column.
$ws.Range("E2").Value = "More"
$ws.Range("E3").Value = "Less"
$ws.Range("E4").Value = "More"
$ws.Range("E5").Value = "Less"
$ws.Range("E6").Value = "More"
$ws.Range("E7").Value = "More"
$ws.Range("E8").Value = "Less"

# Match the saved selection state from the diff (active cell on the new column header).
$ws.Range("E1").Select() | Out-Null
